# Generate Report for handoff
# Update the "Latest Handoff Datetime" (column D, row 5) for the
# 978a7524-7d3a-4a07-9038-42538f5b36e1 file on both the zh-cn and de-de
# localization status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-25 10:37:39"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-25 10:37:48"
